$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 287; this shifts the existing rows 287-363
# down to 289-365, matching the diff (dimension grows from R363 to R365).
$ws.Rows.Item(287).Resize(2).Insert()

# New row 287 (weekly update - new price observation)
$ws.Cells.Item(287, 1).Value = 6
$ws.Cells.Item(287, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(287, 3).Value = "Metropolitana"
$ws.Cells.Item(287, 4).Value = 44641
$ws.Cells.Item(287, 5).Value = 13
$ws.Cells.Item(287, 6).Value = 100112032
$ws.Cells.Item(287, 7).Value = "Zapallo italiano"
$ws.Cells.Item(287, 8).Value = "Sin especificar"
$ws.Cells.Item(287, 9).Value = "Primera"
$ws.Cells.Item(287, 10).Value = 460
$ws.Cells.Item(287, 11).Value = 12000
$ws.Cells.Item(287, 12).Value = 13000
$ws.Cells.Item(287, 13).Value = 12457
$ws.Cells.Item(287, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(287, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(287, 16).Value = 249
$ws.Cells.Item(287, 17).Value = 50
$ws.Cells.Item(287, 18).Value = "Hortaliza"

# New row 288 (weekly update - new price observation)
$ws.Cells.Item(288, 1).Value = 6
$ws.Cells.Item(288, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(288, 3).Value = "Metropolitana"
$ws.Cells.Item(288, 4).Value = 44641
$ws.Cells.Item(288, 5).Value = 13
$ws.Cells.Item(288, 6).Value = 100112032
$ws.Cells.Item(288, 7).Value = "Zapallo italiano"
$ws.Cells.Item(288, 8).Value = "Sin especificar"
$ws.Cells.Item(288, 9).Value = "Primera"
$ws.Cells.Item(288, 10).Value = 180
$ws.Cells.Item(288, 11).Value = 12000
$ws.Cells.Item(288, 12).Value = 12000
$ws.Cells.Item(288, 13).Value = 12000
$ws.Cells.Item(288, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(288, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(288, 16).Value = 240
$ws.Cells.Item(288, 17).Value = 50
$ws.Cells.Item(288, 18).Value = "Hortaliza"
